$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 970
$ws.Range("B3").Value = 492
$ws.Range("B4").Value = 150
$ws.Range("B5").Value = 45
$ws.Range("B6").Value = 250
